$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.192.18'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'585.16"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').Value = "'173.39"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.46%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  +2.76%  '
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '67.054.83'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('D17').Value = '2.497.79'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('D19').Value = "'10.94"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('D20').Value = "'350.36"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = "'68.84"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('E25').Value = '  +2.22%  '
$ws.Range('D26').Value = "'9.18"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.76%  '
$ws.Range('D27').Value = '2.607.66'
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = "'504.18"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('D31').Value = "'7.74"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.120"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = "'161.85"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.37%  '
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').Value = "'18.19"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('D45').Value = "'143.02"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.72%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0260'
$ws.Range('E46').Value = '  +3.28%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = "'3.48"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').Value = "'0.584"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.49%  '
